$wb = $excel.ActiveWorkbook

# Rename the sheets: drop the leading "FFT_" prefix from each sheet name.
$wb.Worksheets.Item(1).Name = "data_CCM code_FLASH"
$wb.Worksheets.Item(2).Name = "data_CCM code_CCM"
$wb.Worksheets.Item(3).Name = "data_RAM code_FLASH"
$wb.Worksheets.Item(4).Name = "data_RAM code_CCM"

# Round the "intensity" row (row 2) to whole numbers, and the "energy" row
# (row 5) down to a 0-1 scale rounded to 3 decimal places, on every sheet.
foreach ($ws in $wb.Worksheets) {
    foreach ($col in @("B", "C", "D")) {
        $intensityCell = $ws.Range($col + "2")
        $intensityCell.Value = [Math]::Round([double]$intensityCell.Value(), 0)

        $energyCell = $ws.Range($col + "5")
        $energyCell.Value = [Math]::Round([double]$energyCell.Value() / 1000, 3)
    }
}
